$d = $word.ActiveDocument

# 1. Insert a new paragraph right after paragraph 1, carrying the text that
#    used to live in paragraph 2 ("还好，日子一天天过去，时光飞逝，随遇而安。").
#    Anchoring the insertion on paragraph 1's end (rather than paragraph 2's
#    start) means the new paragraph mark inherits paragraph 1's eastAsia
#    font hint, matching the target markup.
$firstPara = $d.Paragraphs.Item(1)
$insertionPoint = $firstPara.Range
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(2)
$newPara.Range.Text = "还好，日子一天天过去，时光飞逝，随遇而安。"

# 2. The old paragraph 2 (now paragraph 3, still holding the _GoBack
#    bookmark) gets new wording for the week's diary entry. Restrict the
#    Find/Replace to that paragraph's own range so the freshly inserted
#    paragraph (which now has identical old text) is left untouched.
$thirdPara = $d.Paragraphs.Item(3)
$targetRange = $thirdPara.Range
$targetRange.Find.Execute("还好，日子一天天过去，时光飞逝，随遇而安。", $true, $false, $false, $false, $false,
                           $true, 1, $false, "今天星期五，一周又过去了，天气阴。", 2)
